# "separate dept from affiliations"
#
# The old "dept" column on "PI hours" actually held each PI's full list of
# affiliations (e.g. "['ME', 'AE', 'CSL']"). This edit:
#   - keeps that list, but moves it to a new "app" column
#   - replaces "dept" with just the PI's primary department
#   - renames "dept hours" -> "department hours" and recomputes it to be
#     hours/percentage per primary department
#   - adds a new "unit(accumulative) hours" sheet that keeps the original
#     (pre-edit) per-unit breakdown, now labelled "unit(accumulative)"

$wb = $excel.ActiveWorkbook

# Step 1: snapshot the current "dept hours" sheet as the new
# "unit(accumulative) hours" sheet, placed right after it, before we touch
# its data.
$wsDept = $wb.Worksheets.Item("dept hours")
$wsDept.Copy($null, $wsDept)
$wsUnit = $wb.Worksheets.Item($wsDept.Index + 1)
$wsUnit.Name = "unit(accumulative) hours"
$wsUnit.Range("B1").Value = "unit(accumulative)"

# Step 2: rename "dept hours" -> "department hours" and rewrite its rows to
# the simplified per-department totals (3 rows instead of 5).
$wsDept.Name = "department hours"
$wsDept.Range("B2").Value = "ME"
$wsDept.Range("C2").Value = 39
$wsDept.Range("D2").Value = 70.90909090909091
$wsDept.Range("B3").Value = "ABE"
$wsDept.Range("C3").Value = 8
$wsDept.Range("D3").Value = 14.54545454545454
$wsDept.Range("B4").Value = "ECE"
$wsDept.Range("C4").Value = 8
$wsDept.Range("D4").Value = 14.54545454545454
$wsDept.Range("A5:D6").Delete()

# Step 3: on "PI hours", split the old "dept" column (E, which held the full
# affiliation list) into a simple "dept" column and a new "app" column that
# keeps the full affiliation list.
$wsPI = $wb.Worksheets.Item("PI hours")

# Move the existing list values from E2:E4 into the new F2:F4 "app" column.
$wsPI.Range("E2:E4").Copy()
$wsPI.Range("F2").PasteSpecial(-4163)

# New header "app" in F1, matching the style of the other header cells.
$wsPI.Range("E1").Copy()
$wsPI.Range("F1").PasteSpecial(-4122)
$wsPI.Range("F1").Value = "app"

# Replace column E with the simplified primary department.
$wsPI.Range("E2").Value = "ME"
$wsPI.Range("E3").Value = "ABE"
$wsPI.Range("E4").Value = "ECE"

$wsPI.Activate()
